$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new BOM row (row 3) for the U.FL-R-SMT-1 (40) RF connector.
$ws.Range("A3").Value = "U.FL-R-SMT-1 (40)"
$ws.Range("B3").Value = 2

# Price is stored as text ("1.06"), matching how it was entered originally.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1.06"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 2044
$ws.Range("E3").Value = "https://www.digikey.ch/product-detail/de/hirose-electric-co-ltd/U-FL-R-SMT-1-40/H125761CT-ND/8594840"
$ws.Range("F3").Value = 250

# New column header added after filling in the new row.
$ws.Range("F1").Value = "max reflow temp (°C)"

# Leave the selection where the author last clicked.
$ws.Range("E2").Select()
